$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "MIE(7.891175580313517, 4.0541935852593, -9.40002195553489, 10.556955510359565)"
$ws.Range("C2").Value = "NIG(1.0431245840843828, 0.7791191676431977, 4.50706208858047, 5.968798929786272)"
$ws.Range("D2").Value = "JSU(-0.8133277738103957, 1.0059530038116793, 1.0845497313291559, 2.4060132146351934)"
$ws.Range("E2").Value = "NCT(3.547639280763107, 1.92462760791359, -0.016949325035053066, 4.384565370800585)"

$wb.Save()
